$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in SUZANNE's hours (column I) for week 1 tasks (rows 3-8)
$ws.Range("I3").Value = 0.041666666666666664
$ws.Range("I4").Value = 0.16666666666666666
$ws.Range("I5").Value = 0.041666666666666664
$ws.Range("I6").Value = 0.08333333333333333
$ws.Range("I7").Value = 0.08333333333333333
$ws.Range("I8").Value = 0.041666666666666664

# I4 gets the same time-of-day number format used by the other "hh:mm:ss" cells (style 13)
$ws.Range("I4").NumberFormat = $ws.Range("I12").NumberFormat

# Fill in SUZANNE's hours (column I) for week 2 tasks (rows 11-15)
$ws.Range("I11").Value = 0
$ws.Range("I11").NumberFormat = $ws.Range("I12").NumberFormat

$ws.Range("I13").Value = 0.013888888888888888
$ws.Range("I13").NumberFormat = $ws.Range("I12").NumberFormat

$ws.Range("I14").Value = 0.25

# Totaal tijd besteed week 2 (row 16) for SUZANNE, summing I11:I15
$ws.Range("I16").Formula = "=SUM(I11,I12,I13,I14,I15)"
$ws.Range("I16").NumberFormat = "[h]:mm:ss"

# Totaal tijd project (row 18) for SUZANNE now sums the two week totals instead of week2 detail rows
$ws.Range("I18").Formula = "=SUM(I9,I16)"

# Move the active selection as in the authored workbook
$ws.Range("H14").Select()
